$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data cells to match the latest scrape.
# Cells whose new value is a plain numeric-looking string are first marked as
# Text ("@") format so Excel keeps the exact text (with trailing zeros, etc.)
# instead of silently re-interpreting it as a floating point number.

$ws.Range("D2").Value = '63.112.83'
$ws.Range("E2").Value = '  +5.50%  '
$ws.Range("D3").Value = '2.439.80'
$ws.Range("E3").Value = '  +5.85%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.89'
$ws.Range("E5").Value = '  +4.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.06'
$ws.Range("E6").Value = '  +9.05%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  +2.73%  '
$ws.Range("D9").Value = '2.439.43'
$ws.Range("E9").Value = '  +5.89%  '
$ws.Range("E10").Value = '  +4.10%  '
$ws.Range("E11").Value = '  +3.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.151'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("E13").Value = '  +5.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.27'
$ws.Range("E14").Value = '  +13.10%  '
$ws.Range("D15").Value = '2.879.59'
$ws.Range("E15").Value = '  +6.01%  '
$ws.Range("D16").Value = '62.987.95'
$ws.Range("E16").Value = '  +5.29%  '
$ws.Range("E17").Value = '  +8.55%  '
$ws.Range("D18").Value = '2.451.30'
$ws.Range("E18").Value = '  +6.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.25'
$ws.Range("E19").Value = '  +7.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '340.34'
$ws.Range("E20").Value = '  +9.81%  '
$ws.Range("E21").Value = '  +4.70%  '
$ws.Range("E22").Value = '  +4.58%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.64'
$ws.Range("E24").Value = '  -0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.62'
$ws.Range("E25").Value = '  +3.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.175'
$ws.Range("E26").Value = '  +3.57%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.54'
$ws.Range("E28").Value = '  +15.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.28'
$ws.Range("E29").Value = '  +7.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.35'
$ws.Range("E30").Value = '  +13.03%  '
$ws.Range("D31").Value = '0.0₃0797'
$ws.Range("E31").Value = '  +10.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.83'
$ws.Range("E32").Value = '  +6.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.56'
$ws.Range("E33").Value = '  +13.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '174.65'
$ws.Range("E34").Value = '  +1.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.48'
$ws.Range("E35").Value = '  +10.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.400'
$ws.Range("E36").Value = '  +5.97%  '
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.68'
$ws.Range("E37").Value = '  +5.94%  '
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '375.53'
$ws.Range("E38").Value = '  +18.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.48'
$ws.Range("E39").Value = '  +12.16%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.71'
$ws.Range("E42").Value = '  +13.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.90'
$ws.Range("E43").Value = '  +6.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '145.52'
$ws.Range("E44").Value = '  +6.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.68'
$ws.Range("E45").Value = '  +7.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '20.55'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.594'
$ws.Range("E47").Value = '  +4.90%  '
$ws.Range("E48").Value = '  +6.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0951'
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("E50").Value = '  +5.23%  '
$ws.Range("E51").Value = '  +7.35%  '
